$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tr = $sh.TextFrame2.TextRange
$para = $tr.Paragraphs(1, 1)

# Justify the paragraph
$para.ParagraphFormat.Alignment = 4  # ppAlignJustify

# Rewrite the three runs to fix the wording error
$para.Runs(1, 1).Text = "House renting is a web application that is people can advertise their own house and "
$para.Runs(2, 1).Text = "some can "
$para.Runs(3, 1).Text = "view house advertisement who want to rent house."
